$d = $word.ActiveDocument

# The paragraph "opdracht Groothandel 2:" currently consists of two runs
# ("opdracht" and " Groothandel 2:") separated by grammar-check
# <w:proofErr> markers. The edit merges them into a single run with the
# full text and drops the now-stale proofErr markers.
$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w14:paraId="139F9A24" w14:textId="2F68803F" w:rsidR="007F3972" w:rsidRDefault="0C37DD1C" w:rsidP="39112147">
            <w:pPr>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000" w:themeColor="text1"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
                <w:color w:val="000000" w:themeColor="text1"/>
              </w:rPr>
              <w:t>opdracht Groothandel 2:</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*opdracht*Groothandel 2*") {
        [void]$p.Range.InsertXML($xmlSnippet)
        break
    }
}
